$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.326.94'
$ws.Range("E2").Value = '  +2.46%  '
$ws.Range("D3").Value = '2.107.90'
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '345.38'
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5248'
$ws.Range("E7").Value = '  +2.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4442'
$ws.Range("E8").Value = '  +1.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.72'
$ws.Range("E9").Value = '  +4.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09486'
$ws.Range("E10").Value = '  +4.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.174'
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.04'
$ws.Range("E12").Value = '  +1.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.759'
$ws.Range("E13").Value = '  +8.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.927'
$ws.Range("E14").Value = '  +2.84%  '
$ws.Range("D15").Value = '2.077.03'
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '101.83'
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001168'
$ws.Range("E17").Value = '  +2.25%  '
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("E19").Value = '  +1.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06728'
$ws.Range("E20").Value = '  +1.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.306'
$ws.Range("E21").Value = '  +2.53%  '
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D23").Value = '30.376.32'
$ws.Range("E23").Value = '  +2.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.64'
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.313'
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("D26").Value = '2.339.35'
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.05'
$ws.Range("E27").Value = '  +1.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '163.92'
$ws.Range("E28").Value = '  +0.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.540'
$ws.Range("E29").Value = '  +1.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.55'
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("E31").Value = '  +2.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.739'
$ws.Range("E32").Value = '  +7.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.1056'
$ws.Range("E33").Value = '  +1.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.917'
$ws.Range("E34").Value = '  +15.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.268'
$ws.Range("E35").Value = '  +2.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.924'
$ws.Range("E36").Value = '  -1.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.49'
$ws.Range("E37").Value = '  +3.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02630'
$ws.Range("E38").Value = '  +2.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06811'
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.7036'
$ws.Range("E40").Value = '  +3.11%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.349'
$ws.Range("E41").Value = '  +5.50%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.59'
$ws.Range("E42").Value = '  +2.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.2230'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6844'
$ws.Range("E44").Value = '  +2.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.55'
$ws.Range("E45").Value = '  +3.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.358'
$ws.Range("E46").Value = '  +3.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.356'
$ws.Range("E48").Value = '  +15.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.647'
$ws.Range("E49").Value = '  +1.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000345'
$ws.Range("E50").Value = '  +2.70%  '
$ws.Range("E51").Value = '  +0.85%  '
